# Apply "update scripts wuth new tpm" edit:
#  - recompute TPM-derived metrics for the remaining sending-cluster rows
#  - drop the "Inflammatory-Mac" and "Resolving-Mac" target-cluster rows
#    (the "MuSCs" row shifts up to become row 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 6 (Resolving-Mac) and row 4 (Inflammatory-Mac) first, from
# bottom to top so row indices of rows we still need stay valid. This
# leaves the former row 5 (MuSCs) as the new row 4.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()

# Row 2 (ECs) - new TPM-derived values
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.579684
$ws.Range("H2").Value = 1.739052
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4581623333333333
$ws.Range("N2").Value = 1.374487
$ws.Range("O2").Value = 0.1060599910922654
$ws.Range("P2").Value = 0.1060599910922654
$ws.Range("Q2").Value = 0.265589374036
$ws.Range("R2").Value = 2.390304366324
$ws.Range("S2").Value = 0.1060599910922654
$ws.Range("T2").Value = 0.1060599910922654

# Row 3 (FAPs) - new TPM-derived values
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.579684
$ws.Range("H3").Value = 1.739052
$ws.Range("M3").Value = 0.075101
$ws.Range("O3").Value = 0.01738512926863672
$ws.Range("P3").Value = 0.01738512926863672
$ws.Range("Q3").Value = 0.043534848084
$ws.Range("R3").Value = 0.391813632756
$ws.Range("S3").Value = 0.01738512926863672
$ws.Range("T3").Value = 0.01738512926863672

# Row 4 (was row 5, MuSCs) - new TPM-derived values
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.579684
$ws.Range("H4").Value = 1.739052
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.786578
$ws.Range("N4").Value = 11.359734
$ws.Range("O4").Value = 0.8765548796390978
$ws.Range("P4").Value = 0.8765548796390978
$ws.Range("Q4").Value = 2.195018681352
$ws.Range("R4").Value = 19.755168132168
$ws.Range("S4").Value = 0.8765548796390978
$ws.Range("T4").Value = 0.8765548796390978
